$wb = $excel.ActiveWorkbook

# --- Sheet1: APPInfoDataReqDTO (columns A-I) ---
$ws1 = $wb.Worksheets.Item("APPInfoDataReqDTO")
$ws1.Cells.Item(188, 1).Value = "2022-07-04 21:08:40"
$ws1.Cells.Item(188, 3).Value = "10F872226797"
$ws1.Cells.Item(188, 4).Value = 1656961723
$ws1.Cells.Item(188, 5).Value = 3
$ws1.Cells.Item(188, 6).Value = 12
$ws1.Cells.Item(188, 7).Value = 1
$ws1.Cells.Item(188, 8).Value = 0
$ws1.Cells.Item(188, 9).Value = 1

$ws1.Cells.Item(189, 1).Value = "2022-07-05 10:47:52"
$ws1.Cells.Item(189, 3).Value = "10F872226797"
$ws1.Cells.Item(189, 4).Value = 1657010878
$ws1.Cells.Item(189, 5).Value = 3
$ws1.Cells.Item(189, 6).Value = 12
$ws1.Cells.Item(189, 7).Value = 1
$ws1.Cells.Item(189, 8).Value = 0
$ws1.Cells.Item(189, 9).Value = 1

$ws1.Cells.Item(190, 1).Value = "2022-07-05 11:38:52"
$ws1.Cells.Item(190, 3).Value = "10F872226797"
$ws1.Cells.Item(190, 4).Value = 1657013938
$ws1.Cells.Item(190, 5).Value = 3
$ws1.Cells.Item(190, 6).Value = 12
$ws1.Cells.Item(190, 7).Value = 1
$ws1.Cells.Item(190, 8).Value = 0
$ws1.Cells.Item(190, 9).Value = 1

$ws1.Cells.Item(191, 1).Value = "2022-07-05 18:08:15"
$ws1.Cells.Item(191, 3).Value = "10F872226797"
$ws1.Cells.Item(191, 4).Value = 1657037303
$ws1.Cells.Item(191, 5).Value = 3
$ws1.Cells.Item(191, 6).Value = 12
$ws1.Cells.Item(191, 7).Value = 1
$ws1.Cells.Item(191, 8).Value = 0
$ws1.Cells.Item(191, 9).Value = 1

# --- Sheet2: APPDtuInfoMO (columns A-X) ---
$ws2 = $wb.Worksheets.Item("APPDtuInfoMO")
$ws2.Cells.Item(188, 1).Value = "2022-07-04 21:08:40"
$ws2.Cells.Item(188, 3).Value = 0
$ws2.Cells.Item(188, 4).Value = 522
$ws2.Cells.Item(188, 5).Value = 37122
$ws2.Cells.Item(188, 6).Value = 0
$ws2.Cells.Item(188, 7).Value = 655618
$ws2.Cells.Item(188, 8).Value = 10485764
$ws2.Cells.Item(188, 9).Value = 1
$ws2.Cells.Item(188, 10).Value = 1656961667
$ws2.Cells.Item(188, 11).Value = 62
$ws2.Cells.Item(188, 13).Value = "1.0.1"
$ws2.Cells.Item(188, 15).Value = 21
$ws2.Cells.Item(188, 16).Value = 0
$ws2.Cells.Item(188, 17).Value = 0
$ws2.Cells.Item(188, 18).Value = 0
$ws2.Cells.Item(188, 19).Value = 0
$ws2.Cells.Item(188, 20).Value = 0
$ws2.Cells.Item(188, 21).Value = 0
$ws2.Cells.Item(188, 22).Value = 0
$ws2.Cells.Item(188, 23).Value = 0

$ws2.Cells.Item(189, 1).Value = "2022-07-05 10:47:52"
$ws2.Cells.Item(189, 3).Value = 0
$ws2.Cells.Item(189, 4).Value = 522
$ws2.Cells.Item(189, 5).Value = 37122
$ws2.Cells.Item(189, 6).Value = 0
$ws2.Cells.Item(189, 7).Value = 655618
$ws2.Cells.Item(189, 8).Value = 10485764
$ws2.Cells.Item(189, 9).Value = 1
$ws2.Cells.Item(189, 10).Value = 1657010831
$ws2.Cells.Item(189, 11).Value = 52
$ws2.Cells.Item(189, 13).Value = "1.0.1"
$ws2.Cells.Item(189, 15).Value = 21
$ws2.Cells.Item(189, 16).Value = 0
$ws2.Cells.Item(189, 17).Value = 0
$ws2.Cells.Item(189, 18).Value = 0
$ws2.Cells.Item(189, 19).Value = 0
$ws2.Cells.Item(189, 20).Value = 0
$ws2.Cells.Item(189, 21).Value = 0
$ws2.Cells.Item(189, 22).Value = 0
$ws2.Cells.Item(189, 23).Value = 0

$ws2.Cells.Item(190, 1).Value = "2022-07-05 11:38:52"
$ws2.Cells.Item(190, 3).Value = 0
$ws2.Cells.Item(190, 4).Value = 522
$ws2.Cells.Item(190, 5).Value = 37122
$ws2.Cells.Item(190, 6).Value = 0
$ws2.Cells.Item(190, 7).Value = 0
$ws2.Cells.Item(190, 8).Value = 0
$ws2.Cells.Item(190, 9).Value = 1
$ws2.Cells.Item(190, 10).Value = 1657013925
$ws2.Cells.Item(190, 11).Value = 58
$ws2.Cells.Item(190, 13).Value = "1.0.1"
$ws2.Cells.Item(190, 15).Value = 21
$ws2.Cells.Item(190, 16).Value = 0
$ws2.Cells.Item(190, 17).Value = 0
$ws2.Cells.Item(190, 18).Value = 0
$ws2.Cells.Item(190, 19).Value = 0
$ws2.Cells.Item(190, 20).Value = 0
$ws2.Cells.Item(190, 21).Value = 0
$ws2.Cells.Item(190, 22).Value = 0
$ws2.Cells.Item(190, 23).Value = 0

$ws2.Cells.Item(191, 1).Value = "2022-07-05 18:08:15"
$ws2.Cells.Item(191, 3).Value = 0
$ws2.Cells.Item(191, 4).Value = 522
$ws2.Cells.Item(191, 5).Value = 37122
$ws2.Cells.Item(191, 6).Value = 0
$ws2.Cells.Item(191, 7).Value = 655618
$ws2.Cells.Item(191, 8).Value = 10485764
$ws2.Cells.Item(191, 9).Value = 1
$ws2.Cells.Item(191, 10).Value = 1657037247
$ws2.Cells.Item(191, 11).Value = 52
$ws2.Cells.Item(191, 13).Value = "1.0.1"
$ws2.Cells.Item(191, 15).Value = 21
$ws2.Cells.Item(191, 16).Value = 0
$ws2.Cells.Item(191, 17).Value = 0
$ws2.Cells.Item(191, 18).Value = 0
$ws2.Cells.Item(191, 19).Value = 0
$ws2.Cells.Item(191, 20).Value = 0
$ws2.Cells.Item(191, 21).Value = 0
$ws2.Cells.Item(191, 22).Value = 0
$ws2.Cells.Item(191, 23).Value = 0

# --- Sheet3: APPPvInfoMO (columns A-M) ---
$ws3 = $wb.Worksheets.Item("APPPvInfoMO")
$ws3.Cells.Item(560, 1).Value = "2022-07-04 21:08:40"
$ws3.Cells.Item(560, 3).Value = 0
$ws3.Cells.Item(560, 4).Value = 19110224228993
$ws3.Cells.Item(560, 5).Value = 0
$ws3.Cells.Item(560, 6).Value = 10012
$ws3.Cells.Item(560, 7).Value = 269627400
$ws3.Cells.Item(560, 8).Value = 256
$ws3.Cells.Item(560, 9).Value = 2560
$ws3.Cells.Item(560, 10).Value = 8193
$ws3.Cells.Item(560, 11).Value = 0
$ws3.Cells.Item(560, 12).Value = 0
$ws3.Cells.Item(560, 13).Value = 0

$ws3.Cells.Item(561, 1).Value = "2022-07-04 21:08:40"
$ws3.Cells.Item(561, 3).Value = 0
$ws3.Cells.Item(561, 4).Value = 19110224226576
$ws3.Cells.Item(561, 5).Value = 0
$ws3.Cells.Item(561, 6).Value = 10012
$ws3.Cells.Item(561, 7).Value = 269627400
$ws3.Cells.Item(561, 8).Value = 256
$ws3.Cells.Item(561, 9).Value = 2560
$ws3.Cells.Item(561, 10).Value = 8193
$ws3.Cells.Item(561, 11).Value = 0
$ws3.Cells.Item(561, 12).Value = 0
$ws3.Cells.Item(561, 13).Value = 0

$ws3.Cells.Item(562, 1).Value = "2022-07-04 21:08:40"
$ws3.Cells.Item(562, 3).Value = 0
$ws3.Cells.Item(562, 4).Value = 19110224228646
$ws3.Cells.Item(562, 5).Value = 0
$ws3.Cells.Item(562, 6).Value = 10012
$ws3.Cells.Item(562, 7).Value = 269627400
$ws3.Cells.Item(562, 8).Value = 256
$ws3.Cells.Item(562, 9).Value = 2560
$ws3.Cells.Item(562, 10).Value = 8193
$ws3.Cells.Item(562, 11).Value = 0
$ws3.Cells.Item(562, 12).Value = 0
$ws3.Cells.Item(562, 13).Value = 0

$ws3.Cells.Item(563, 1).Value = "2022-07-05 10:47:52"
$ws3.Cells.Item(563, 3).Value = 0
$ws3.Cells.Item(563, 4).Value = 19110224228993
$ws3.Cells.Item(563, 5).Value = 0
$ws3.Cells.Item(563, 6).Value = 10012
$ws3.Cells.Item(563, 7).Value = 269627400
$ws3.Cells.Item(563, 8).Value = 256
$ws3.Cells.Item(563, 9).Value = 2560
$ws3.Cells.Item(563, 10).Value = 8193
$ws3.Cells.Item(563, 11).Value = 0
$ws3.Cells.Item(563, 12).Value = 0
$ws3.Cells.Item(563, 13).Value = 0

$ws3.Cells.Item(564, 1).Value = "2022-07-05 10:47:52"
$ws3.Cells.Item(564, 3).Value = 0
$ws3.Cells.Item(564, 4).Value = 19110224226576
$ws3.Cells.Item(564, 5).Value = 0
$ws3.Cells.Item(564, 6).Value = 10012
$ws3.Cells.Item(564, 7).Value = 269627400
$ws3.Cells.Item(564, 8).Value = 256
$ws3.Cells.Item(564, 9).Value = 2560
$ws3.Cells.Item(564, 10).Value = 8193
$ws3.Cells.Item(564, 11).Value = 0
$ws3.Cells.Item(564, 12).Value = 0
$ws3.Cells.Item(564, 13).Value = 0

$ws3.Cells.Item(565, 1).Value = "2022-07-05 10:47:52"
$ws3.Cells.Item(565, 3).Value = 0
$ws3.Cells.Item(565, 4).Value = 19110224228646
$ws3.Cells.Item(565, 5).Value = 0
$ws3.Cells.Item(565, 6).Value = 10012
$ws3.Cells.Item(565, 7).Value = 269627400
$ws3.Cells.Item(565, 8).Value = 256
$ws3.Cells.Item(565, 9).Value = 2560
$ws3.Cells.Item(565, 10).Value = 8193
$ws3.Cells.Item(565, 11).Value = 0
$ws3.Cells.Item(565, 12).Value = 0
$ws3.Cells.Item(565, 13).Value = 0

$ws3.Cells.Item(566, 1).Value = "2022-07-05 11:38:52"
$ws3.Cells.Item(566, 3).Value = 0
$ws3.Cells.Item(566, 4).Value = 19110224228993
$ws3.Cells.Item(566, 5).Value = 0
$ws3.Cells.Item(566, 6).Value = 10012
$ws3.Cells.Item(566, 7).Value = 269627400
$ws3.Cells.Item(566, 8).Value = 256
$ws3.Cells.Item(566, 9).Value = 2560
$ws3.Cells.Item(566, 10).Value = 8193
$ws3.Cells.Item(566, 11).Value = 0
$ws3.Cells.Item(566, 12).Value = 0
$ws3.Cells.Item(566, 13).Value = 0

$ws3.Cells.Item(567, 1).Value = "2022-07-05 11:38:52"
$ws3.Cells.Item(567, 3).Value = 0
$ws3.Cells.Item(567, 4).Value = 19110224226576
$ws3.Cells.Item(567, 5).Value = 0
$ws3.Cells.Item(567, 6).Value = 10012
$ws3.Cells.Item(567, 7).Value = 269627400
$ws3.Cells.Item(567, 8).Value = 256
$ws3.Cells.Item(567, 9).Value = 2560
$ws3.Cells.Item(567, 10).Value = 8193
$ws3.Cells.Item(567, 11).Value = 0
$ws3.Cells.Item(567, 12).Value = 0
$ws3.Cells.Item(567, 13).Value = 0

$ws3.Cells.Item(568, 1).Value = "2022-07-05 11:38:52"
$ws3.Cells.Item(568, 3).Value = 0
$ws3.Cells.Item(568, 4).Value = 19110224228646
$ws3.Cells.Item(568, 5).Value = 0
$ws3.Cells.Item(568, 6).Value = 10012
$ws3.Cells.Item(568, 7).Value = 269627400
$ws3.Cells.Item(568, 8).Value = 256
$ws3.Cells.Item(568, 9).Value = 2560
$ws3.Cells.Item(568, 10).Value = 8193
$ws3.Cells.Item(568, 11).Value = 0
$ws3.Cells.Item(568, 12).Value = 0
$ws3.Cells.Item(568, 13).Value = 0

$ws3.Cells.Item(569, 1).Value = "2022-07-05 18:08:15"
$ws3.Cells.Item(569, 3).Value = 0
$ws3.Cells.Item(569, 4).Value = 19110224228993
$ws3.Cells.Item(569, 5).Value = 0
$ws3.Cells.Item(569, 6).Value = 10012
$ws3.Cells.Item(569, 7).Value = 269627400
$ws3.Cells.Item(569, 8).Value = 256
$ws3.Cells.Item(569, 9).Value = 2560
$ws3.Cells.Item(569, 10).Value = 8193
$ws3.Cells.Item(569, 11).Value = 0
$ws3.Cells.Item(569, 12).Value = 0
$ws3.Cells.Item(569, 13).Value = 0

$ws3.Cells.Item(570, 1).Value = "2022-07-05 18:08:15"
$ws3.Cells.Item(570, 3).Value = 0
$ws3.Cells.Item(570, 4).Value = 19110224226576
$ws3.Cells.Item(570, 5).Value = 0
$ws3.Cells.Item(570, 6).Value = 10012
$ws3.Cells.Item(570, 7).Value = 269627400
$ws3.Cells.Item(570, 8).Value = 256
$ws3.Cells.Item(570, 9).Value = 2560
$ws3.Cells.Item(570, 10).Value = 8193
$ws3.Cells.Item(570, 11).Value = 0
$ws3.Cells.Item(570, 12).Value = 0
$ws3.Cells.Item(570, 13).Value = 0

$ws3.Cells.Item(571, 1).Value = "2022-07-05 18:08:15"
$ws3.Cells.Item(571, 3).Value = 0
$ws3.Cells.Item(571, 4).Value = 19110224228646
$ws3.Cells.Item(571, 5).Value = 0
$ws3.Cells.Item(571, 6).Value = 10012
$ws3.Cells.Item(571, 7).Value = 269627400
$ws3.Cells.Item(571, 8).Value = 256
$ws3.Cells.Item(571, 9).Value = 2560
$ws3.Cells.Item(571, 10).Value = 8193
$ws3.Cells.Item(571, 11).Value = 0
$ws3.Cells.Item(571, 12).Value = 0
$ws3.Cells.Item(571, 13).Value = 0

